$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $addr, $val)
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws 'D2' '62.974.63'
Set-TextCell $ws 'E2' '  -5.08%  '

# Row 3
Set-TextCell $ws 'D3' '3.093.92'
Set-TextCell $ws 'E3' '  -6.40%  '

# Row 5
Set-TextCell $ws 'D5' '555.85'
Set-TextCell $ws 'E5' '  -5.62%  '

# Row 6
Set-TextCell $ws 'D6' '159.30'
Set-TextCell $ws 'E6' '  -11.29%  '

# Row 7
Set-TextCell $ws 'E7' '  +0.00%  '

# Row 8
Set-TextCell $ws 'D8' '0.576'
Set-TextCell $ws 'E8' '  -10.15%  '

# Row 9
Set-TextCell $ws 'D9' '3.095.47'
Set-TextCell $ws 'E9' '  -6.21%  '

# Row 10
Set-TextCell $ws 'D10' '6.68'
Set-TextCell $ws 'E10' '  -2.29%  '

# Row 11
Set-TextCell $ws 'D11' '0.113'
Set-TextCell $ws 'E11' '  -9.72%  '

# Row 12
Set-TextCell $ws 'D12' '0.372'
Set-TextCell $ws 'E12' '  -7.19%  '

# Row 13
Set-TextCell $ws 'D13' '3.639.68'
Set-TextCell $ws 'E13' '  -6.09%  '

# Row 14
Set-TextCell $ws 'E14' '  -1.74%  '

# Row 15
Set-TextCell $ws 'D15' '63.011.59'
Set-TextCell $ws 'E15' '  -5.04%  '

# Row 16
Set-TextCell $ws 'D16' '24.28'
Set-TextCell $ws 'E16' '  -8.92%  '

# Row 17
Set-TextCell $ws 'D17' '3.096.21'
Set-TextCell $ws 'E17' '  -5.90%  '

# Row 18
Set-TextCell $ws 'D18' '0.0000150'
Set-TextCell $ws 'E18' '  -7.82%  '

# Row 19
Set-TextCell $ws 'D19' '392.76'
Set-TextCell $ws 'E19' '  -7.80%  '

# Row 20
Set-TextCell $ws 'D20' '12.25'
Set-TextCell $ws 'E20' '  -6.03%  '

# Row 21
Set-TextCell $ws 'D21' '5.12'
Set-TextCell $ws 'E21' '  -6.75%  '

# Row 22
Set-TextCell $ws 'D22' '6.92'
Set-TextCell $ws 'E22' '  -5.29%  '

# Row 23
Set-TextCell $ws 'D23' '0.999'
Set-TextCell $ws 'E23' '  -0.14%  '

# Row 24
Set-TextCell $ws 'E24' '  +0.26%  '

# Row 25
Set-TextCell $ws 'D25' '66.77'
Set-TextCell $ws 'E25' '  -6.18%  '

# Row 26
Set-TextCell $ws 'E26' '  -4.15%  '

# Row 27
Set-TextCell $ws 'D27' '0.470'
Set-TextCell $ws 'E27' '  -8.16%  '

# Row 28
Set-TextCell $ws 'D28' '0.0₃0986'
Set-TextCell $ws 'E28' '  -13.57%  '

# Row 29
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextCell $ws 'D29' '1.00'
Set-TextCell $ws 'E29' '  +0.16%  '

# Row 30
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell $ws 'D30' '8.49'
Set-TextCell $ws 'E30' '  -9.45%  '

# Row 31
Set-TextCell $ws 'E31' '  -0.08%  '

# Row 32
Set-TextCell $ws 'D32' '1.76'
Set-TextCell $ws 'E32' '  -8.20%  '

# Row 33
Set-TextCell $ws 'D33' '20.61'
Set-TextCell $ws 'E33' '  -7.65%  '

# Row 34
Set-TextCell $ws 'D34' '4.74'
Set-TextCell $ws 'E34' '  -8.43%  '

# Row 35
Set-TextCell $ws 'D35' '6.13'
Set-TextCell $ws 'E35' '  -6.71%  '

# Row 36
Set-TextCell $ws 'D36' '1.08'
Set-TextCell $ws 'E36' '  -9.26%  '

# Row 37
Set-TextCell $ws 'D37' '150.41'
Set-TextCell $ws 'E37' '  -5.41%  '

# Row 38
Set-TextCell $ws 'D38' '1.29'
Set-TextCell $ws 'E38' '  -10.07%  '

# Row 39
Set-TextCell $ws 'D39' '2.671.82'
Set-TextCell $ws 'E39' '  -6.55%  '

# Row 40
Set-TextCell $ws 'D40' '1.62'
Set-TextCell $ws 'E40' '  -10.12%  '

# Row 41
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell $ws 'D41' '22.98'
Set-TextCell $ws 'E41' '  -12.24%  '

# Row 42
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell $ws 'D42' '3.99'
Set-TextCell $ws 'E42' '  -7.84%  '

# Row 43
Set-TextCell $ws 'D43' '38.11'
Set-TextCell $ws 'E43' '  -3.80%  '

# Row 44
Set-TextCell $ws 'D44' '0.690'
Set-TextCell $ws 'E44' '  -7.74%  '

# Row 45
Set-TextCell $ws 'D45' '0.0599'
Set-TextCell $ws 'E45' '  -6.15%  '

# Row 46
Set-TextCell $ws 'D46' '5.38'
Set-TextCell $ws 'E46' '  -8.58%  '

# Row 47
Set-TextCell $ws 'D47' '0.0251'
Set-TextCell $ws 'E47' '  -7.01%  '

# Row 48
$ws.Range('B48').Value = 'FirstDigitalUSD'
$ws.Range('C48').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextCell $ws 'D48' '0.999'
Set-TextCell $ws 'E48' '  +0.08%  '

# Row 49
$ws.Range('B49').Value = 'Bittensor'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell $ws 'D49' '279.70'
Set-TextCell $ws 'E49' '  -10.81%  '

# Row 50
Set-TextCell $ws 'D50' '20.33'
Set-TextCell $ws 'E50' '  -11.24%  '

# Row 51
Set-TextCell $ws 'D51' '0.0965'
Set-TextCell $ws 'E51' '  -5.63%  '

Write-Output "Update complete"